$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 359, shifting existing rows 359:495 down to 360:496.
$ws.Rows("359:359").Insert()

# Populate the newly inserted row with the new data record.
$ws.Range("A359").Value = 4
$ws.Range("B359").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C359").Value = "Los Lagos"
$ws.Range("D359").Value = 45146
$ws.Range("E359").Value = 10
$ws.Range("F359").Value = 100112037
$ws.Range("G359").Value = "Cebollín"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 180
$ws.Range("K359").Value = 5500
$ws.Range("L359").Value = 6000
$ws.Range("M359").Value = 5750
$ws.Range("N359").Value = "`$/paquete 36 unidades"
$ws.Range("O359").Value = "Región Metropolitana"
$ws.Range("P359").Value = 160
$ws.Range("Q359").Value = 36
$ws.Range("R359").Value = "Hortaliza"
